# [Fonds de solidarite] Add 2022-06-23 data
# Update nombre_aides (column C) and montant_total (column E) for the rows
# affected by the refreshed publication date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 13;  C = 187863;  E = 1168289738 },
    @{ Row = 88;  C = 71283;   E = 110329263 },
    @{ Row = 91;  C = 18884;   E = 75370883 },
    @{ Row = 93;  C = 16953;   E = 50862471 },
    @{ Row = 98;  C = 6301;    E = 19502085 },
    @{ Row = 100; C = 9348;    E = 23899308 },
    @{ Row = 121; C = 1306416; E = 2275552863 },
    @{ Row = 122; C = 382;     E = 1260995 },
    @{ Row = 129; C = 633827;  E = 3435576311 },
    @{ Row = 130; C = 4250;    E = 141666559 },
    @{ Row = 132; C = 586042;  E = 3473327690 },
    @{ Row = 136; C = 26706;   E = 144413770 },
    @{ Row = 178; C = 515892;  E = 891219559 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
